# Auto-generated Excel COM-interop script to apply scheduled price-refresh update
# to the Omega_Profits workbook (per diff: updates currentAveragePrice* / Leve*Price / Leve*Profit
# columns, H-N, across 8 item-crafting sheets).

$wb = $excel.ActiveWorkbook


# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 11: 'Gotta Bounce' / 'Rubber' (item id 5533)
$ws.Range("H11").Value = 184.5238
$ws.Range("I11").Value = 184.5238
$ws.Range("K11").Value = 184.5238
$ws.Range("M11").Value = -44.52379999999999
# Row 68: "Can't Sleep, Inquisitors Will Eat Me" / 'Wyvernskin Grimoire' (item id 10647)
$ws.Range("H68").Value = 169998
$ws.Range("J68").Value = 169998
$ws.Range("L68").Value = 169998
$ws.Range("N68").Value = -171496
# Row 71: 'Allow No Fallacies (L)' / 'Wyvernskin Grimoire' (item id 10647)
$ws.Range("H71").Value = 169998
$ws.Range("J71").Value = 169998
$ws.Range("L71").Value = 509994
$ws.Range("N71").Value = -517482
# Row 127: 'Liquid Competence' / "Competent Craftsman's Draught" (item id 36114)
$ws.Range("H127").Value = 3064.3333
$ws.Range("I127").Value = 2098
$ws.Range("J127").Value = 4997
$ws.Range("K127").Value = 6294
$ws.Range("L127").Value = 14991
$ws.Range("M127").Value = -1334
$ws.Range("N127").Value = -24911
# Row 129: 'Practical Command' / "Commanding Craftsman's Draught" (item id 36115)
$ws.Range("H129").Value = 2555
$ws.Range("J129").Value = 2600
$ws.Range("L129").Value = 7800
$ws.Range("N129").Value = -17800
# Row 132: 'Fast-forwarding Flora' / 'Growth Formula Lambda' (item id 44049)
$ws.Range("H132").Value = 2435.3125
$ws.Range("I132").Value = 2090
$ws.Range("K132").Value = 6270
$ws.Range("M132").Value = -3740

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2: "Ain't Got No Ingots" / 'Bronze Ingot' (item id 27713)
$ws.Range("H2").Value = 4673.1816
$ws.Range("J2").Value = 7006.5
$ws.Range("L2").Value = 7006.5
$ws.Range("N2").Value = -7232.5
# Row 19: 'Stadium Envy' / 'Bronze Gauntlets' (item id 3550)
$ws.Range("H19").Value = 7004.25
$ws.Range("I19").Value = 8336
$ws.Range("K19").Value = 8336
$ws.Range("M19").Value = -8107
# Row 45: 'Hollow Hallmarks' / 'Mythril Ingot' (item id 27714)
$ws.Range("H45").Value = 14702.625
$ws.Range("I45").Value = 21696.2
$ws.Range("K45").Value = 21696.2
$ws.Range("M45").Value = -21319.2
# Row 61: 'Dealing with the Tough Stuff' / 'Cobalt Ingot' (item id 43999)
$ws.Range("H61").Value = 4210.7295
$ws.Range("I61").Value = 4079.9143
$ws.Range("K61").Value = 4079.9143
$ws.Range("M61").Value = -3867.9143
# Row 76: 'Sometimes the South Wins' / 'Titanium Mail of Fending' (item id 10679)
$ws.Range("H76").Value = 72933
$ws.Range("J76").Value = 72933
$ws.Range("L76").Value = 72933
$ws.Range("N76").Value = -73609
# Row 79: 'The Thriller of Autumn (L)' / 'Titanium Mail of Fending' (item id 10679)
$ws.Range("H79").Value = 72933
$ws.Range("J79").Value = 72933
$ws.Range("L79").Value = 72933
$ws.Range("N79").Value = -75273
# Row 110: 'Scheduled Maintenance' / 'Deepgold Ingot' (item id 27708)
$ws.Range("H110").Value = 1451.4
$ws.Range("I110").Value = 1734
$ws.Range("J110").Value = 886.2
$ws.Range("K110").Value = 1734
$ws.Range("L110").Value = 886.2
$ws.Range("M110").Value = 311
$ws.Range("N110").Value = -4976.2
# Row 116: 'No Scope' / 'Titanbronze Ingot' (item id 27713)
$ws.Range("H116").Value = 4673.1816
$ws.Range("J116").Value = 7006.5
$ws.Range("L116").Value = 7006.5
$ws.Range("N116").Value = -11594.5
# Row 122: 'Haste for High Durium' / 'High Durium Nugget' (item id 36168)
$ws.Range("H122").Value = 1926.0555
$ws.Range("I122").Value = 1947.875
$ws.Range("K122").Value = 5843.625
$ws.Range("M122").Value = -3393.625
# Row 132: "Don't Bore Me, Ore Me" / 'Mountain Chromite Ingot' (item id 43997)
$ws.Range("H132").Value = 3373.8667
$ws.Range("I132").Value = 3382.75
$ws.Range("K132").Value = 10148.25
$ws.Range("M132").Value = -7618.25
# Row 136: 'Metal with Mettle' / 'Cobalt Tungsten Ingot' (item id 43999)
$ws.Range("H136").Value = 4210.7295
$ws.Range("I136").Value = 4079.9143
$ws.Range("K136").Value = 12239.7429
$ws.Range("M136").Value = -9689.742899999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3: 'Hells Bells' / 'Bronze Ingot' (item id 27713)
$ws.Range("H3").Value = 4673.1816
$ws.Range("J3").Value = 7006.5
$ws.Range("L3").Value = 7006.5
$ws.Range("N3").Value = -7234.5
# Row 86: 'Through Thick and Thin' / 'Adamantite Nugget' (item id 12526)
$ws.Range("H86").Value = 4492.533
$ws.Range("I86").Value = 3239.5
$ws.Range("J86").Value = 6998.6
$ws.Range("K86").Value = 3239.5
$ws.Range("L86").Value = 6998.6
$ws.Range("M86").Value = -2116.5
$ws.Range("N86").Value = -9244.6
# Row 89: 'Piercing Eyes Deserve Piercing Shafts (L)' / 'Adamantite Nugget' (item id 12526)
$ws.Range("H89").Value = 4492.533
$ws.Range("I89").Value = 3239.5
$ws.Range("J89").Value = 6998.6
$ws.Range("K89").Value = 16197.5
$ws.Range("L89").Value = 34993
$ws.Range("M89").Value = -10581.5
$ws.Range("N89").Value = -46225
# Row 92: 'Have Blade, Will Travel' / 'High Steel Katzbalger' (item id 18033)
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").ClearContents()
$ws.Range("N92").Value = 0

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 10: 'Spears and Sorcery' / 'Maple Crook' (item id 1997)
$ws.Range("H10").Value = 992
$ws.Range("I10").Value = 138
$ws.Range("K10").Value = 138
$ws.Range("M10").Value = 1
# Row 31: 'Wall Not Found' / 'Walnut Lumber' (item id 44023)
$ws.Range("H31").Value = 2845.0942
$ws.Range("I31").Value = 2090
$ws.Range("J31").Value = 4091
$ws.Range("K31").Value = 2090
$ws.Range("L31").Value = 4091
$ws.Range("M31").Value = -1795
$ws.Range("N31").Value = -4681
# Row 34: 'Armoires of the Rich and Famous' / 'Walnut Lumber' (item id 44023)
$ws.Range("H34").Value = 2845.0942
$ws.Range("I34").Value = 2090
$ws.Range("J34").Value = 4091
$ws.Range("K34").Value = 2090
$ws.Range("L34").Value = 4091
$ws.Range("M34").Value = -1888
$ws.Range("N34").Value = -4495
# Row 86: 'Birch, Please' / 'Birch Lumber' (item id 12584)
$ws.Range("H86").Value = 11910300
$ws.Range("I86").Value = 13893979
$ws.Range("K86").Value = 13893979
$ws.Range("M86").Value = -13892856
# Row 89: 'Built This City on Blocks and Soul (L)' / 'Birch Lumber' (item id 12584)
$ws.Range("H89").Value = 11910300
$ws.Range("I89").Value = 13893979
$ws.Range("K89").Value = 69469895
$ws.Range("M89").Value = -69464279
# Row 103: 'Spare a Rod and Spoil the Fishers' / 'Gazelle Horn Fishing Rod' (item id 19558)
$ws.Range("H103").Value = 59666
$ws.Range("I103").Value = 54499.5
$ws.Range("J103").Value = 69999
$ws.Range("K103").Value = 54499.5
$ws.Range("L103").Value = 69999
$ws.Range("M103").Value = -53327.5
$ws.Range("N103").Value = -72343

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 113: "Can't Eat Just One" / 'Night Vinegar' (item id 27843)
$ws.Range("H113").Value = 2289.2727
$ws.Range("J113").Value = 3021.4614
$ws.Range("L113").Value = 9064.3842
$ws.Range("N113").Value = -13404.3842

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 102: 'Put the Metal to the Peddle' / 'Durium Ingot' (item id 36169)
$ws.Range("H102").Value = 1947.25
$ws.Range("I102").Value = 1967.3478
$ws.Range("K102").Value = 1967.3478
$ws.Range("M102").Value = -345.3478
# Row 126: 'Gold Rush Order' / 'Phrygian Gold Ingot' (item id 36184)
$ws.Range("H126").Value = 6236.909
$ws.Range("I126").Value = 3766.1667
$ws.Range("J126").Value = 9201.799999999999
$ws.Range("K126").Value = 11298.5001
$ws.Range("L126").Value = 27605.4
$ws.Range("M126").Value = -8828.500100000001
$ws.Range("N126").Value = -32545.4
# Row 132: 'On Board for Lar' / 'Lar Ingot' (item id 44008)
$ws.Range("H132").Value = 2347.3462
$ws.Range("I132").Value = 1732.5883
$ws.Range("J132").Value = 3508.5557
$ws.Range("K132").Value = 5197.7649
$ws.Range("L132").Value = 10525.6671
$ws.Range("M132").Value = -2667.7649
$ws.Range("N132").Value = -15585.6671

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 14: 'Quelling Bloody Rumors' / 'Hard Leather Shoes' (item id 3771)
$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").ClearContents()
$ws.Range("N14").Value = 0
# Row 16: 'Saddle Sore' / 'Hard Leather' (item id 5289)
$ws.Range("H16").Value = 3450
$ws.Range("I16").Value = 3549.8
$ws.Range("K16").Value = 3549.8
$ws.Range("M16").Value = -3379.8
# Row 125: 'Scouting Talent' / 'Luncheon Toadskin Jacket of Scouting' (item id 34271)
$ws.Range("H125").Value = 82997.39999999999
$ws.Range("J125").Value = 82997.39999999999
$ws.Range("L125").Value = 82997.39999999999
$ws.Range("N125").Value = -92837.39999999999
# Row 132: 'Tenets of Tanning' / 'Silver Lobo Leather' (item id 44058)
$ws.Range("H132").Value = 27567.275
$ws.Range("I132").Value = 40001.156
$ws.Range("K132").Value = 120003.468
$ws.Range("M132").Value = -117473.468

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 4: 'Not Cool Enough' / 'Hempen Undershirt' (item id 2996)
$ws.Range("H4").Value = 8584.5
$ws.Range("I4").Value = 16494
$ws.Range("K4").Value = 16494
$ws.Range("M4").Value = -16381
# Row 14: 'Hat in Hand' / 'Straw Hat' (item id 2658)
$ws.Range("H14").Value = 1037.9762
$ws.Range("I14").Value = 1047.375
$ws.Range("J14").Value = 850
$ws.Range("K14").Value = 1047.375
$ws.Range("L14").Value = 850
$ws.Range("M14").Value = -879.375
$ws.Range("N14").Value = -1186
# Row 122: 'Heavy Armoire' / 'Dark Hempen Cloth' (item id 36208)
$ws.Range("H122").Value = 3635.442
$ws.Range("I122").Value = 2748.6177
$ws.Range("K122").Value = 8245.8531
$ws.Range("M122").Value = -5795.8531
# Row 126: 'A Polished Purchase' / 'Snow Linen' (item id 36210)
$ws.Range("H126").Value = 6156
$ws.Range("I126").Value = 5663.8
$ws.Range("K126").Value = 16991.4
$ws.Range("M126").Value = -14521.4
# Row 132: 'Comfy Cabins' / 'Snow Cotton Cloth' (item id 44029)
$ws.Range("H132").Value = 1922
$ws.Range("J132").Value = 2076.6
$ws.Range("L132").Value = 6229.799999999999
$ws.Range("N132").Value = -11289.8
# Row 136: 'Weaving the Envelope' / 'Sarcenet Cloth' (item id 44031)
$ws.Range("H136").Value = 1587.1904
$ws.Range("I136").Value = 1531.9
$ws.Range("K136").Value = 4595.700000000001
$ws.Range("M136").Value = -2045.700000000001

Write-Host "Applied scheduled price refresh."
